$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menus_and_forms")

# Replace the "sheet_name" header in B1 with "menu_or_form"
$ws.Range("B1").Value = "menu_or_form"

# Update the active cell selection from B4 to B2
$ws.Range("B2").Select()
